# Weekly CompStat refresh: roll the report forward by one week
# (Volume Number 16 -> 17, week-of dates +7 days) and load the newly
# collected crime-complaint figures for the 105th Precinct table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Masthead text: "Volume 32   Number  16" -> "...17"
#    and the "Report Covering the Week  4/14/2025  Through  4/20/2025"
#    line -> "...4/21/2025  Through  4/27/2025". Edit only the runs that
#    change (via Characters, like retyping just that portion in Excel)
#    so the rest of each string is left untouched.
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "17"
$ws.Range("C9").Characters(27, 9).Text = "4/21/2025"
$ws.Range("C9").Characters(47, 9).Text = "4/27/2025"

# ---------------------------------------------------------------------
# 2. Helper: convert a cell that currently holds a plain number into a
#    text cell holding one of the sheet's existing shared label strings
#    ("0" or "***.*"), without disturbing its visual style. We copy the
#    number format from a donor cell that already carries that exact
#    style, then copy the value from a donor cell that already carries
#    that exact text, so no new style/shared-string entries are minted.
# ---------------------------------------------------------------------
function Set-LabelCell($target, $donor) {
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($target).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($target).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Donors already on the sheet (row 14, untouched by this edit):
#   C14 = s13/text "0"      E14 = s13/text "***.*"      I14 = s14/number 1
$ZeroLabel = "C14"
$DashLabel = "E14"
$OneNumber = "I14"

# Rape (row 15): 28-day 2024 count drops to 0 -> "%Chg" becomes "***.*"
Set-LabelCell "G15" $ZeroLabel
Set-LabelCell "H15" $DashLabel

# Robbery (row 16): WTD 2024 count drops to 0 -> "%Chg" becomes "***.*"
Set-LabelCell "D16" $ZeroLabel
Set-LabelCell "E16" $DashLabel

# Burglary (row 18): WTD 2025 count drops to 0
Set-LabelCell "C18" $ZeroLabel

# UCR Rape* (row 27): 28-day 2024 count drops to 0 -> "%Chg" becomes "***.*"
Set-LabelCell "G27" $ZeroLabel
Set-LabelCell "H27" $DashLabel

# Other Sex Crimes (row 28): both WTD counts drop to 0 -> "%Chg" becomes "***.*"
Set-LabelCell "C28" $ZeroLabel
Set-LabelCell "D28" $ZeroLabel
Set-LabelCell "E28" $DashLabel

# Hate Crimes (row 31): WTD 2025 and 28-day 2025 counts go from "0" to 1
Set-LabelCell "C31" $OneNumber
Set-LabelCell "F31" $OneNumber

# ---------------------------------------------------------------------
# 3. Plain numeric refreshes across the crime-complaint table
#    (new weekly totals and their recomputed percentages).
# ---------------------------------------------------------------------

# Rape (row 15)
$ws.Range("N15").Value = -78.571428571428

# Robbery (row 16)
$ws.Range("C16").Value = 4
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 250
$ws.Range("I16").Value = 30
$ws.Range("K16").Value = 100
$ws.Range("L16").Value = 11.111111111111
$ws.Range("M16").Value = -70.588235294117
$ws.Range("N16").Value = -90.654205607476

# Fel. Assault (row 17)
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -10
$ws.Range("I17").Value = 67
$ws.Range("J17").Value = 92
$ws.Range("K17").Value = -27.173913043478
$ws.Range("L17").Value = -2.898550724637
$ws.Range("M17").Value = -27.173913043478
$ws.Range("N17").Value = -33.663366336633

# Burglary (row 18)
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = 21.428571428571
$ws.Range("L18").Value = -1.923076923076
$ws.Range("M18").Value = -55.263157894736
$ws.Range("N18").Value = -90

# Gr. Larceny (row 19)
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 20
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 2.941176470588
$ws.Range("I19").Value = 125
$ws.Range("J19").Value = 134
$ws.Range("K19").Value = -6.716417910447
$ws.Range("L19").Value = 8.695652173913
$ws.Range("M19").Value = -16.666666666666
$ws.Range("N19").Value = -26.035502958579

# G.L.A. (row 20)
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -42.857142857142
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 58
$ws.Range("J20").Value = 76
$ws.Range("K20").Value = -23.684210526315
$ws.Range("L20").Value = 5.454545454545
$ws.Range("M20").Value = -50.847457627118
$ws.Range("N20").Value = -94.722474977252

# TOTAL (row 21)
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -15.384615384615
$ws.Range("F21").Value = 86
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = -10.416666666666
$ws.Range("I21").Value = 335
$ws.Range("J21").Value = 365
$ws.Range("K21").Value = -8.219178082191
$ws.Range("L21").Value = 4.361370716510
$ws.Range("M21").Value = -43.027210884353
$ws.Range("N21").Value = -84.909909909909

# Petit Larceny (row 24)
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = -54.545454545454
$ws.Range("F24").Value = 26
$ws.Range("G24").Value = 44
$ws.Range("H24").Value = -40.909090909090
$ws.Range("I24").Value = 163
$ws.Range("J24").Value = 218
$ws.Range("K24").Value = -25.229357798165
$ws.Range("L24").Value = -34.8
$ws.Range("M24").Value = -36.078431372549

# Retail Theft (row 25)
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = -50
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 28.571428571428
$ws.Range("I25").Value = 38
$ws.Range("J25").Value = 49
$ws.Range("K25").Value = -22.448979591836
$ws.Range("L25").Value = -13.636363636363

# Misd. Assault (row 26)
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -9.090909090909
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = 91.304347826087
$ws.Range("I26").Value = 144
$ws.Range("J26").Value = 106
$ws.Range("K26").Value = 35.849056603773
$ws.Range("L26").Value = 37.142857142857
$ws.Range("M26").Value = -27.638190954773

# Hate Crimes (row 31)
$ws.Range("I31").Value = 2
$ws.Range("K31").Value = 100
$ws.Range("L31").Value = 100
